# Update column F (dSF) values for specific rows to match repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = -1
    "F8"  = -2
    "F14" = -6
    "F17" = -7
    "F25" = -4
    "F27" = -4
    "F30" = -5
    "F33" = -9
    "F35" = -3
    "F36" = 3
    "F37" = -3
    "F40" = -3
    "F44" = -4
    "F47" = 4
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
